$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9445658922195435
$ws.Range("B1").Value = 2.14838981628418
$ws.Range("C1").Value = 8.517941474914551
$ws.Range("D1").Value = 1.763637185096741
$ws.Range("E1").Value = 1.404836416244507
